$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B and C column updates (coin name / link swap for rows 37-38)
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

# D column (Price) updates - force text to avoid numeric coercion, then clear the
# number-format style iron_native applies so cells keep their original (unstyled) look
$ws.Range("D2").Value = '''60.775.03'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = '''3.373.48'
$ws.Range("D3").ClearFormats()
$ws.Range("D5").Value = '''568.18'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = '''135.88'
$ws.Range("D6").ClearFormats()
$ws.Range("D8").Value = '''3.374.89'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = '''0.467'
$ws.Range("D9").ClearFormats()
$ws.Range("D13").Value = '''3.944.82'
$ws.Range("D13").ClearFormats()
$ws.Range("D15").Value = '''26.03'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = '''3.371.59'
$ws.Range("D16").ClearFormats()
$ws.Range("D18").Value = '''60.893.00'
$ws.Range("D18").ClearFormats()
$ws.Range("D20").Value = '''13.61'
$ws.Range("D20").ClearFormats()
$ws.Range("D22").Value = '''371.34'
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = '''3.508.33'
$ws.Range("D23").ClearFormats()
$ws.Range("D25").Value = '''1.00'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = '''70.79'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = '''0.0000122'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").Value = '''0.176'
$ws.Range("D28").ClearFormats()
$ws.Range("D31").Value = '''7.29'
$ws.Range("D31").ClearFormats()
$ws.Range("D35").Value = '''23.24'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").Value = '''5.10'
$ws.Range("D36").ClearFormats()
$ws.Range("D37").Value = '''1.53'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").Value = '''6.76'
$ws.Range("D38").ClearFormats()
$ws.Range("D39").Value = '''164.68'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").Value = '''0.0757'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = '''41.66'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = '''24.99'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = '''1.69'
$ws.Range("D45").ClearFormats()
$ws.Range("D48").Value = '''2.523.44'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = '''23.44'
$ws.Range("D49").ClearFormats()
$ws.Range("D51").Value = '''2.40'
$ws.Range("D51").ClearFormats()

# E column (Volume 1h %) updates
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("E11").Value = '  -3.95%  '
$ws.Range("E12").Value = '  -2.77%  '
$ws.Range("E13").Value = '  -0.71%  '
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("E17").Value = '  -4.56%  '
$ws.Range("E18").Value = '  -1.38%  '
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("E20").Value = '  -4.58%  '
$ws.Range("E21").Value = '  -2.84%  '
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("E23").Value = '  -0.49%  '
$ws.Range("E24").Value = '  -2.38%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("E27").Value = '  -3.73%  '
$ws.Range("E28").Value = '  +9.34%  '
$ws.Range("E29").Value = '  -3.89%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  -3.48%  '
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("E36").Value = '  -4.84%  '
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("E47").Value = '  -6.53%  '
$ws.Range("E48").Value = '  +7.31%  '
$ws.Range("E49").Value = '  +2.88%  '
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("E51").Value = '  +2.84%  '
